# Remove logging system configuration
#
# Appends one additional sensor-log data row (row 34) to each of the four
# worksheets. Each new row duplicates the preceding (last) logged row,
# but with the timestamp advanced by one hour - mirroring the existing
# hourly cadence of the log.

$wb = $excel.ActiveWorkbook

$newRows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = "2025-03-05 17:42:06"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x d";  F = 400; G = "568631262647113770877196"; H = 400; I = 13  },
    @{ Sheet = "ROW35-MID-LIFTER"; A = "2025-03-05 17:29:35"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x90,"; E = "0x e";  F = 400; G = "568631262647113770942732"; H = 400; I = 14  },
    @{ Sheet = "ROW02-FE-LIFTER";  A = "2025-03-05 17:51:45"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; D = "0x01,0x90,"; E = "0xff";  F = 400; G = "568631262647113769959692"; H = 400; I = 255 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = "2025-03-05 17:41:15"; B = "0x01,0x90 "; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x90,"; E = "0x 3";  F = 400; G = "568631262647113769959692"; H = 400; I = 3   }
)

foreach ($entry in $newRows) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $r = 34

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F

    # Column G holds a long, purely-numeric identifier that must remain
    # a text string (it would otherwise overflow double precision and
    # be mangled into scientific notation), so force a text format
    # before assigning it.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $entry.G

    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}
